# CH-77 Character-Based Rhombus — "Created a generalized odd solution"
#
# 1. Refactor the existing fixed 13x13 odd-rhombus array formula (C67:O79) to
#    use LET with zero-based r/c helper variables (rz/cz) instead of inlining
#    "r-1"/"c-1" everywhere. Same result, clearer formula.
# 2. Add an explanatory note in V65.
# 3. Add a new, generalized (parameterized by N) odd-rhombus array formula in
#    V67, spilling into V67:AF77 (an 11x11 example), to show the general LET
#    formula that will work for any odd N.
# 4. Narrow the columns under the new array (V:AF) to match the other grid
#    columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDA")
$ws.Activate() | Out-Null

# --- 1. Rewrite the original 13x13 rhombus formula with rz/cz via LET -----
$ws.Range("C67:O79").FormulaArray = "=MAKEARRAY(13,13,LAMBDA(r,c,LET(rz,r-1,cz,c-1,IF(AND(rz-cz<=6,rz+cz>=6,rz+cz<=18,rz-cz>=-6),""*"",""""))))"

# --- 2. Explanatory note cell ---------------------------------------------
$ws.Range("V65").Value = "This solves the odd case. The even case will be similar and won't teach me anything new."

# --- 3. Generalized odd-rhombus formula, demoed at N=11 -------------------
$ws.Range("V67").Formula2 = "=LET(N,11,MAKEARRAY(N,N,LAMBDA(r,c,LET(rz,r-1,cz,c-1,nz,INT(N/2),IF(AND(rz-cz<=nz,rz+cz>=nz,rz+cz<=N*1.5-1.5,rz-cz>=-nz),""*"","""")))))"

# --- 4. Column widths for the new grid (V:AF) ------------------------------
$ws.Range("V1:AF1").ColumnWidth = 1.7271205357142856

# --- Selection / view, matches where the author left off -------------------
$ws.Range("V65").Select() | Out-Null
